$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 167.70732
$ws.Range("I12").Value = 157.44
$ws.Range("J12").Value = 183.75
$ws.Range("K12").Value = 157.44
$ws.Range("L12").Value = 183.75
$ws.Range("M12").Value = 12.56
$ws.Range("N12").Value = -523.75
$ws.Range("H17").Value = 1263215.2
$ws.Range("J17").Value = 1437960.2
$ws.Range("L17").Value = 4313880.6
$ws.Range("N17").Value = -4314216.6
$ws.Range("H19").Value = 5411
$ws.Range("I19").Value = 946.3333
$ws.Range("J19").Value = 8759.5
$ws.Range("K19").Value = 946.3333
$ws.Range("L19").Value = 8759.5
$ws.Range("M19").Value = -771.3333
$ws.Range("N19").Value = -9109.5
$ws.Range("H28").Value = 260189
$ws.Range("I28").Value = 335583.34
$ws.Range("K28").Value = 335583.34
$ws.Range("M28").Value = -335098.34
$ws.Range("H33").Value = 14750772
$ws.Range("J33").Value = 41668000
$ws.Range("L33").Value = 41668000
$ws.Range("N33").Value = -41668458
$ws.Range("H40").Value = 2435.7144
$ws.Range("I40").Value = 2116.6667
$ws.Range("K40").Value = 2116.6667
$ws.Range("M40").Value = -1941.6667
$ws.Range("H64").Value = 5166.6665
$ws.Range("I64").Value = 4000
$ws.Range("J64").Value = 5750
$ws.Range("K64").Value = 4000
$ws.Range("L64").Value = 5750
$ws.Range("M64").Value = -3752
$ws.Range("N64").Value = -6246
$ws.Range("H67").Value = 5166.6665
$ws.Range("I67").Value = 4000
$ws.Range("J67").Value = 5750
$ws.Range("K67").Value = 4000
$ws.Range("L67").Value = 5750
$ws.Range("M67").Value = -3142
$ws.Range("N67").Value = -7466
$ws.Range("H76").Value = 14058219
$ws.Range("I76").Value = 3994.75
$ws.Range("K76").Value = 3994.75
$ws.Range("M76").Value = -3679.75
$ws.Range("H79").Value = 14058219
$ws.Range("I79").Value = 3994.75
$ws.Range("K79").Value = 3994.75
$ws.Range("M79").Value = -2902.75
$ws.Range("H80").Value = 1202956.1
$ws.Range("I80").Value = 3804878.5
$ws.Range("J80").Value = 2068.8462
$ws.Range("K80").Value = 11414635.5
$ws.Range("L80").Value = 6206.5386
$ws.Range("M80").Value = -11413637.5
$ws.Range("N80").Value = -8202.5386
$ws.Range("H83").Value = 1202956.1
$ws.Range("I83").Value = 3804878.5
$ws.Range("J83").Value = 2068.8462
$ws.Range("K83").Value = 34243906.5
$ws.Range("L83").Value = 18619.6158
$ws.Range("M83").Value = -34238914.5
$ws.Range("N83").Value = -28603.6158
$ws.Range("H98").Value = 1319.625
$ws.Range("I98").Value = 1065.8518
$ws.Range("K98").Value = 1065.8518
$ws.Range("M98").Value = 432.1482000000001
$ws.Range("H111").Value = 8631.75
$ws.Range("I111").Value = 2827.5
$ws.Range("K111").Value = 8482.5
$ws.Range("M111").Value = -5415.5
$ws.Range("H113").Value = 166673170
$ws.Range("I113").Value = 500004000
$ws.Range("J113").Value = 7747.5
$ws.Range("K113").Value = 500004000
$ws.Range("L113").Value = 7747.5
$ws.Range("M113").Value = -500000746
$ws.Range("N113").Value = -14255.5
$ws.Range("H116").Value = 41748892
$ws.Range("I116").Value = 25113116
$ws.Range("J116").Value = 83338340
$ws.Range("K116").Value = 25113116
$ws.Range("L116").Value = 83338340
$ws.Range("M116").Value = -25109674
$ws.Range("N116").Value = -83345224
$ws.Range("H118").Value = 345
$ws.Range("I118").Value = 345
$ws.Range("K118").Value = 1035
$ws.Range("M118").Value = 622
$ws.Range("H122").Value = 1319.625
$ws.Range("I122").Value = 1065.8518
$ws.Range("K122").Value = 3197.5554
$ws.Range("M122").Value = -747.5553999999997
$ws.Range("H129").Value = 1451.8096
$ws.Range("I129").Value = 1101.4667
$ws.Range("J129").Value = 2327.6667
$ws.Range("K129").Value = 3304.4001
$ws.Range("L129").Value = 6983.000100000001
$ws.Range("M129").Value = 1695.5999
$ws.Range("N129").Value = -16983.0001
$ws.Range("H132").Value = 3888.5
$ws.Range("I132").Value = 3485.551
$ws.Range("K132").Value = 10456.653
$ws.Range("M132").Value = -7926.653
$ws.Range("H135").Value = 47620836
$ws.Range("I135").Value = 71429330
$ws.Range("J135").Value = 3848.2856
$ws.Range("K135").Value = 642863970
$ws.Range("L135").Value = 34634.5704
$ws.Range("M135").Value = -642861435
$ws.Range("N135").Value = -39704.5704
$ws.Range("H137").Value = 2078.018
$ws.Range("I137").Value = 1789.5128
$ws.Range("J137").Value = 2781.25
$ws.Range("K137").Value = 5368.538399999999
$ws.Range("L137").Value = 8343.75
$ws.Range("M137").Value = -2818.538399999999
$ws.Range("N137").Value = -13443.75
$ws.Range("H138").Value = 2007.6
$ws.Range("I138").Value = 1036.5
$ws.Range("J138").Value = 4504.7144
$ws.Range("K138").Value = 3109.5
$ws.Range("L138").Value = 13514.1432
$ws.Range("M138").Value = 2030.5
$ws.Range("N138").Value = -23794.1432
$ws.Range("H141").Value = 1637.8214
$ws.Range("I141").Value = 1034.6
$ws.Range("J141").Value = 6664.6665
$ws.Range("K141").Value = 3103.8
$ws.Range("L141").Value = 19993.9995
$ws.Range("M141").Value = 2076.2
$ws.Range("N141").Value = -30353.9995

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1931.4
$ws.Range("I2").Value = 1931.4
$ws.Range("K2").Value = 1931.4
$ws.Range("M2").Value = -1818.4
$ws.Range("H39").Value = 10999.5
$ws.Range("I39").Value = 10999.5
$ws.Range("K39").Value = 10999.5
$ws.Range("M39").Value = -10479.5
$ws.Range("H61").Value = 18183526
$ws.Range("I61").Value = 20409816
$ws.Range("K61").Value = 20409816
$ws.Range("M61").Value = -20409604
$ws.Range("H74").Value = 3080.0715
$ws.Range("I74").Value = 3036.8647
$ws.Range("J74").Value = 3399.8
$ws.Range("K74").Value = 3036.8647
$ws.Range("L74").Value = 3399.8
$ws.Range("M74").Value = -2162.8647
$ws.Range("N74").Value = -5147.8
$ws.Range("H77").Value = 3080.0715
$ws.Range("I77").Value = 3036.8647
$ws.Range("J77").Value = 3399.8
$ws.Range("K77").Value = 15184.3235
$ws.Range("L77").Value = 16999
$ws.Range("M77").Value = -10816.3235
$ws.Range("N77").Value = -25735
$ws.Range("H94").Value = 26475.834
$ws.Range("J94").Value = 26311
$ws.Range("L94").Value = 26311
$ws.Range("N94").Value = -28113
$ws.Range("H104").Value = 26916.166
$ws.Range("J104").Value = 26916.166
$ws.Range("L104").Value = 26916.166
$ws.Range("N104").Value = -33904.166
$ws.Range("H110").Value = 6027
$ws.Range("I110").Value = 6189.684
$ws.Range("K110").Value = 6189.684
$ws.Range("M110").Value = -4144.684
$ws.Range("H116").Value = 1931.4
$ws.Range("I116").Value = 1931.4
$ws.Range("K116").Value = 1931.4
$ws.Range("M116").Value = 362.5999999999999
$ws.Range("H119").Value = 44250
$ws.Range("J119").Value = 44250
$ws.Range("L119").Value = 44250
$ws.Range("N119").Value = -53926
$ws.Range("H122").Value = 14708777
$ws.Range("I122").Value = 22729606
$ws.Range("K122").Value = 68188818
$ws.Range("M122").Value = -68186368
$ws.Range("H132").Value = 2895.6333
$ws.Range("I132").Value = 2721.8635
$ws.Range("J132").Value = 3373.5
$ws.Range("K132").Value = 8165.5905
$ws.Range("L132").Value = 10120.5
$ws.Range("M132").Value = -5635.5905
$ws.Range("N132").Value = -15180.5
$ws.Range("H136").Value = 18183526
$ws.Range("I136").Value = 20409816
$ws.Range("K136").Value = 61229448
$ws.Range("M136").Value = -61226898

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1931.4
$ws.Range("I3").Value = 1931.4
$ws.Range("K3").Value = 1931.4
$ws.Range("M3").Value = -1817.4
$ws.Range("H20").Value = 22254.883
$ws.Range("I20").Value = 19740
$ws.Range("K20").Value = 19740
$ws.Range("M20").Value = -19493
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H94").Value = 2546.75
$ws.Range("I94").Value = 2706.2
$ws.Range("K94").Value = 2706.2
$ws.Range("M94").Value = -2255.2
$ws.Range("H107").Value = 962.7143
$ws.Range("I107").Value = 962.7143
$ws.Range("K107").Value = 962.7143
$ws.Range("M107").Value = 957.2857
$ws.Range("H134").Value = 1947.9
$ws.Range("I134").Value = 1600.5454
$ws.Range("J134").Value = 2372.4443
$ws.Range("K134").Value = 4801.6362
$ws.Range("L134").Value = 7117.3329
$ws.Range("M134").Value = -2266.6362
$ws.Range("N134").Value = -12187.3329

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2354.6428
$ws.Range("J16").Value = 3000
$ws.Range("L16").Value = 3000
$ws.Range("N16").Value = -3574
$ws.Range("H31").Value = 4947.2837
$ws.Range("I31").Value = 15106.111
$ws.Range("J31").Value = 3370.9138
$ws.Range("K31").Value = 15106.111
$ws.Range("L31").Value = 3370.9138
$ws.Range("M31").Value = -14811.111
$ws.Range("N31").Value = -3960.9138
$ws.Range("H32").Value = 1515.6666
$ws.Range("I32").Value = 299
$ws.Range("J32").Value = 2732.3333
$ws.Range("K32").Value = 299
$ws.Range("L32").Value = 2732.3333
$ws.Range("M32").Value = 17
$ws.Range("N32").Value = -3364.3333
$ws.Range("H34").Value = 4947.2837
$ws.Range("I34").Value = 15106.111
$ws.Range("J34").Value = 3370.9138
$ws.Range("K34").Value = 15106.111
$ws.Range("L34").Value = 3370.9138
$ws.Range("M34").Value = -14904.111
$ws.Range("N34").Value = -3774.9138
$ws.Range("H58").Value = 418611.53
$ws.Range("I58").Value = 1205.5
$ws.Range("J58").Value = 1253423.6
$ws.Range("K58").Value = 1205.5
$ws.Range("L58").Value = 1253423.6
$ws.Range("M58").Value = -1002.5
$ws.Range("N58").Value = -1253829.6
$ws.Range("H62").Value = 6874.6665
$ws.Range("I62").Value = 6590.5454
$ws.Range("K62").Value = 6590.5454
$ws.Range("M62").Value = -5966.5454
$ws.Range("H65").Value = 6874.6665
$ws.Range("I65").Value = 6590.5454
$ws.Range("K65").Value = 32952.727
$ws.Range("M65").Value = -29832.727
$ws.Range("H70").Value = 30059.2
$ws.Range("J70").Value = 30059.2
$ws.Range("L70").Value = 30059.2
$ws.Range("N70").Value = -30689.2
$ws.Range("H73").Value = 30059.2
$ws.Range("J73").Value = 30059.2
$ws.Range("L73").Value = 30059.2
$ws.Range("N73").Value = -32243.2
$ws.Range("H105").Value = 2521.75
$ws.Range("I105").Value = 1359.8
$ws.Range("J105").Value = 4458.3335
$ws.Range("K105").Value = 1359.8
$ws.Range("L105").Value = 4458.3335
$ws.Range("M105").Value = 387.2
$ws.Range("N105").Value = -7952.3335
$ws.Range("H107").Value = 6549.846
$ws.Range("I107").Value = 5187
$ws.Range("K107").Value = 5187
$ws.Range("M107").Value = -3267
$ws.Range("H113").Value = 2354.6428
$ws.Range("J113").Value = 3000
$ws.Range("L113").Value = 3000
$ws.Range("N113").Value = -7340
$ws.Range("H125").Value = 46166.668
$ws.Range("J125").Value = 46166.668
$ws.Range("L125").Value = 46166.668
$ws.Range("N125").Value = -51086.668
$ws.Range("H134").Value = 904.94446
$ws.Range("I134").Value = 726.21875
$ws.Range("J134").Value = 2334.75
$ws.Range("K134").Value = 2178.65625
$ws.Range("L134").Value = 7004.25
$ws.Range("M134").Value = 356.34375
$ws.Range("N134").Value = -12074.25
$ws.Range("H136").Value = 418611.53
$ws.Range("I136").Value = 1205.5
$ws.Range("J136").Value = 1253423.6
$ws.Range("K136").Value = 3616.5
$ws.Range("L136").Value = 3760270.8
$ws.Range("M136").Value = -1066.5
$ws.Range("N136").Value = -3765370.8

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 14850901
$ws.Range("I4").Value = 27842668
$ws.Range("J4").Value = 3166.7856
$ws.Range("K4").Value = 83528004
$ws.Range("L4").Value = 9500.356800000001
$ws.Range("M4").Value = -83527892
$ws.Range("N4").Value = -9724.356800000001
$ws.Range("H31").Value = 2910
$ws.Range("J31").Value = 2910
$ws.Range("L31").Value = 8730
$ws.Range("N31").Value = -9306
$ws.Range("H44").Value = 4499.684
$ws.Range("I44").Value = 3642.5
$ws.Range("K44").Value = 10927.5
$ws.Range("M44").Value = -10529.5
$ws.Range("H45").Value = 1299
$ws.Range("I45").Value = 1299
$ws.Range("K45").Value = 3897
$ws.Range("M45").Value = -3365
$ws.Range("H55").Value = 799.8
$ws.Range("J55").Value = 1000
$ws.Range("L55").Value = 3000
$ws.Range("N55").Value = -3354
$ws.Range("H80").Value = 6238.8696
$ws.Range("I80").Value = 5812.4375
$ws.Range("J80").Value = 7213.5713
$ws.Range("K80").Value = 17437.3125
$ws.Range("L80").Value = 21640.7139
$ws.Range("M80").Value = -16501.3125
$ws.Range("N80").Value = -23512.7139
$ws.Range("H83").Value = 6238.8696
$ws.Range("I83").Value = 5812.4375
$ws.Range("J83").Value = 7213.5713
$ws.Range("K83").Value = 52311.9375
$ws.Range("L83").Value = 64922.14169999999
$ws.Range("M83").Value = -47631.9375
$ws.Range("N83").Value = -74282.14169999999
$ws.Range("H103").Value = 525.8889
$ws.Range("I103").Value = 190.5
$ws.Range("J103").Value = 1196.6666
$ws.Range("K103").Value = 571.5
$ws.Range("L103").Value = 3589.9998
$ws.Range("M103").Value = 307.5
$ws.Range("N103").Value = -5347.9998
$ws.Range("H108").Value = 1472.1818
$ws.Range("I108").Value = 1519.4
$ws.Range("K108").Value = 4558.200000000001
$ws.Range("M108").Value = -1678.200000000001
$ws.Range("H132").Value = 2133.3333
$ws.Range("I132").Value = 2000
$ws.Range("J132").Value = 2320
$ws.Range("K132").Value = 18000
$ws.Range("L132").Value = 20880
$ws.Range("M132").Value = -15470
$ws.Range("N132").Value = -25940
$ws.Range("H134").Value = 992.8333
$ws.Range("I134").Value = 992.8333
$ws.Range("K134").Value = 2978.4999
$ws.Range("M134").Value = 2091.5001
$ws.Range("H137").Value = 41669148
$ws.Range("J137").Value = 3054.5
$ws.Range("L137").Value = 9163.5
$ws.Range("N137").Value = -19363.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3844.25
$ws.Range("I80").Value = 3762.75
$ws.Range("J80").Value = 3898.5833
$ws.Range("K80").Value = 3762.75
$ws.Range("L80").Value = 3898.5833
$ws.Range("M80").Value = -2764.75
$ws.Range("N80").Value = -5894.5833
$ws.Range("H83").Value = 3844.25
$ws.Range("I83").Value = 3762.75
$ws.Range("J83").Value = 3898.5833
$ws.Range("K83").Value = 18813.75
$ws.Range("L83").Value = 19492.9165
$ws.Range("M83").Value = -13821.75
$ws.Range("N83").Value = -29476.9165
$ws.Range("H113").Value = 8310.5
$ws.Range("I113").Value = 7756.5
$ws.Range("J113").Value = 9972.5
$ws.Range("K113").Value = 7756.5
$ws.Range("L113").Value = 9972.5
$ws.Range("M113").Value = -5586.5
$ws.Range("N113").Value = -14312.5
$ws.Range("H126").Value = 7150.4
$ws.Range("I126").Value = 8252.895
$ws.Range("J126").Value = 3659.1667
$ws.Range("K126").Value = 24758.685
$ws.Range("L126").Value = 10977.5001
$ws.Range("M126").Value = -22288.685
$ws.Range("N126").Value = -15917.5001
$ws.Range("H132").Value = 242474.53
$ws.Range("I132").Value = 316406.03
$ws.Range("K132").Value = 949218.0900000001
$ws.Range("M132").Value = -946688.0900000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 728.0357
$ws.Range("I16").Value = 728.0357
$ws.Range("K16").Value = 728.0357
$ws.Range("M16").Value = -558.0357
$ws.Range("H22").Value = 1189.5454
$ws.Range("I22").Value = 1546.6
$ws.Range("K22").Value = 1546.6
$ws.Range("M22").Value = -1251.6
$ws.Range("H27").Value = 1189.5454
$ws.Range("I27").Value = 1546.6
$ws.Range("K27").Value = 1546.6
$ws.Range("M27").Value = -1439.6
$ws.Range("H40").Value = 2898.75
$ws.Range("I40").Value = 2859.1333
$ws.Range("K40").Value = 2859.1333
$ws.Range("M40").Value = -2723.1333
$ws.Range("H82").Value = 3373.3333
$ws.Range("I82").Value = 3538.0588
$ws.Range("J82").Value = 2673.25
$ws.Range("K82").Value = 3538.0588
$ws.Range("L82").Value = 2673.25
$ws.Range("M82").Value = -3177.0588
$ws.Range("N82").Value = -3395.25
$ws.Range("H85").Value = 3373.3333
$ws.Range("I85").Value = 3538.0588
$ws.Range("J85").Value = 2673.25
$ws.Range("K85").Value = 3538.0588
$ws.Range("L85").Value = 2673.25
$ws.Range("M85").Value = -2290.0588
$ws.Range("N85").Value = -5169.25
$ws.Range("H125").Value = 98426.28999999999
$ws.Range("J125").Value = 98426.28999999999
$ws.Range("L125").Value = 98426.28999999999
$ws.Range("N125").Value = -108266.29
$ws.Range("H132").Value = 3684.0217
$ws.Range("I132").Value = 1711.5483
$ws.Range("J132").Value = 7760.467
$ws.Range("K132").Value = 5134.644899999999
$ws.Range("L132").Value = 23281.401
$ws.Range("M132").Value = -2604.644899999999
$ws.Range("N132").Value = -28341.401
$ws.Range("H136").Value = 1967.386
$ws.Range("I136").Value = 1651.5217
$ws.Range("K136").Value = 4954.5651
$ws.Range("M136").Value = -2404.5651

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 6487.706
$ws.Range("I14").Value = 5623.3335
$ws.Range("K14").Value = 5623.3335
$ws.Range("M14").Value = -5455.3335
$ws.Range("H100").Value = 4636.7646
$ws.Range("I100").Value = 5919.4165
$ws.Range("K100").Value = 11838.833
$ws.Range("M100").Value = -11297.833
$ws.Range("H107").Value = 1233.8889
$ws.Range("I107").Value = 567.3333
$ws.Range("K107").Value = 1701.9999
$ws.Range("M107").Value = 218.0001
$ws.Range("H113").Value = 6751.3
$ws.Range("I113").Value = 8167.923
$ws.Range("J113").Value = 4120.4287
$ws.Range("K113").Value = 24503.769
$ws.Range("L113").Value = 12361.2861
$ws.Range("M113").Value = -22333.769
$ws.Range("N113").Value = -16701.2861
$ws.Range("H122").Value = 1987.0588
$ws.Range("I122").Value = 1918.9166
$ws.Range("J122").Value = 2150.6
$ws.Range("K122").Value = 5756.7498
$ws.Range("L122").Value = 6451.799999999999
$ws.Range("M122").Value = -3306.7498
$ws.Range("N122").Value = -11351.8
$ws.Range("H132").Value = 469280.84
$ws.Range("I132").Value = 692999.5600000001
$ws.Range("J132").Value = 5863.5
$ws.Range("K132").Value = 2078998.68
$ws.Range("L132").Value = 17590.5
$ws.Range("M132").Value = -2076468.68
$ws.Range("N132").Value = -22650.5
$ws.Range("H136").Value = 2732.66
$ws.Range("I136").Value = 2005.8788
$ws.Range("J136").Value = 4143.4707
$ws.Range("K136").Value = 6017.636399999999
$ws.Range("L136").Value = 12430.4121
$ws.Range("M136").Value = -3467.636399999999
$ws.Range("N136").Value = -17530.4121

Write-Host "All edits applied."